$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the image path for the "image" column (P) on the data row (row 2).
# This appends a new shared string and a new cell reference, matching the
# commit "Added paths to images in documents".
$ws.Range("P2").Value = "C:Users/vano/Documents/GitHub/ZPI_VAF/iaff_assistant/images/Student/costs.jpg"

# Reflect the user's navigation to the newly-filled cell.
$ws.Range("P2").Select() | Out-Null
